# Applies the commit "inverse availability, support dates, support results upload"
# to the org_data / lecturer_data workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # org_data
$ws2 = $wb.Worksheets.Item(2)   # lecturer_data

# ---------------------------------------------------------------------------
# 1) Fix the "enthicity" typo -> "ethnicity" in both sheets' header rows.
#    (This is the only place the shared string is used, so the rewrite also
#    reshuffles the shared-strings table exactly like the target diff.)
# ---------------------------------------------------------------------------
$ws1.Range("G1").Value = "ethnicity"
$ws2.Range("F1").Value = "ethnicity"

# ---------------------------------------------------------------------------
# 2) org_data: first_date/second_date/third_date columns (H, J, L) used to
#    store bare day-of-month integers. They now hold full Excel date serials
#    (October 2023) formatted as "d-mmm".
# ---------------------------------------------------------------------------
$dates = @{
    2 = @{ H = 45205; J = 45218; L = 45228 }
    3 = @{ H = 45209; J = 45215; L = 45217 }
    4 = @{ H = 45211; J = 45218; L = 45223 }
    5 = @{ H = 45218; J = 45222; L = 45228 }
    6 = @{ H = 45202; J = 45204; L = 45224 }
    7 = @{ H = 45202; J = 45206; L = 45211 }
    8 = @{ H = 45204; J = 45208; L = 45224 }
}

foreach ($row in $dates.Keys) {
    $cols = $dates[$row]
    foreach ($col in $cols.Keys) {
        $ws1.Range("$col$row").Value = $cols[$col]
    }
}

$ws1.Range("H2:H8").NumberFormat = "d-mmm"
$ws1.Range("J2:J8").NumberFormat = "d-mmm"
$ws1.Range("L2:L8").NumberFormat = "d-mmm"

# ---------------------------------------------------------------------------
# 3) View state: org_data becomes the active/selected sheet & tab, with
#    M13 selected; lecturer_data keeps F13 selected but is no longer active.
# ---------------------------------------------------------------------------
$ws2.Range("F13").Select()
$ws1.Select()
$ws1.Range("M13").Select()
